$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for row 1 ---
$ws.Range("L1").Value = "TJLP"
$ws.Range("M1").Value = "SELIC"
$ws.Range("N1").Value = "CDI"
$ws.Range("O1").Value = "Situacao"

# --- Move the "Situacao" status text from column L to column O for rows 3-6 ---
for ($r = 3; $r -le 6; $r++) {
    $val = $ws.Cells.Item($r, 12).Value   # column L = 12
    $ws.Cells.Item($r, 15).Value = $val   # column O = 15
    $ws.Cells.Item($r, 12).ClearContents()
}

# --- New row 7 data ---
$ws.Range("A7").Value = [DateTime]"2026-02-11"
$ws.Range("A7").NumberFormat = "dd/mm/yyyy"

$ws.Range("D7").Value = 5.183
$ws.Range("D7").NumberFormat = "0.0000"
$ws.Range("E7").Value = 5.1836
$ws.Range("E7").NumberFormat = "0.0000"

$ws.Range("H7").Value = 6.146
$ws.Range("H7").NumberFormat = "0.0000"
$ws.Range("I7").Value = 6.1477
$ws.Range("I7").NumberFormat = "0.0000"

$ws.Range("J7").Value = 6.7076
$ws.Range("J7").NumberFormat = "0.0000"
$ws.Range("K7").Value = 6.7102
$ws.Range("K7").NumberFormat = "0.0000"

$ws.Range("L7").Value = 0.0919
$ws.Range("L7").NumberFormat = "yyyy-mm-dd"
$ws.Range("L7").NumberFormat = "0.0000%"
$ws.Range("M7").Value = 0.15
$ws.Range("M7").NumberFormat = "0.0000%"

$ws.Range("N7").Value = 0.0551310642
$ws.Range("N7").NumberFormat = "0.0000000000"

$ws.Range("O7").Value = "OK 11/02/2026 16:53:39"
